$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.587.86"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "2.990.34"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "382.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.53"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0845"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.67"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "3.455.85"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "2.991.43"
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.975"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").Value = "51.615.78"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +5.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.170"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +19.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.10"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.111"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.56%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.39"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0459"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.53"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("E40").Value = "  -5.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.24"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.283"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +21.62%  "
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.52%  "
$ws.Range("D48").Value = "2.054.32"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.26"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +12.30%  "
$ws.Range("E51").Value = "  +4.47%  "
